# Daily attendance processing - 2026-01-19 15:12:42
# Normalize the "Recorded By" (column G) author lists:
#   - "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   - "backup@backdoor.com, System, system"   -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "backup@backdoor.com, System, system") {
        $cell.Value2 = "backup@backdoor.com, system, System"
    }
}
